$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.382126988660822
$ws1.Range("C2").Value = -0.6681233342238342
$ws1.Range("B3").Value = -0.4033999489948966
$ws1.Range("C3").Value = -1.757645729525741
$ws1.Range("B4").Value = -0.7855548477543319
$ws1.Range("C4").Value = 0.5634447424687831

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -0.5836999949484416
$ws2.Range("C2").Value = -0.306531637635053
$ws2.Range("B3").Value = -0.7818529999851362
$ws2.Range("C3").Value = 0.2187089838408385
$ws2.Range("B4").Value = -1.229024437938767
$ws2.Range("C4").Value = 0.9610437075293694
